$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scan timestamps in column E (rows 2-13) to reflect new scan times
$ws.Range("E2").Value = "2025-04-11 13:21:07"
$ws.Range("E3").Value = "2025-04-11 13:21:11"
$ws.Range("E4").Value = "2025-04-11 13:21:14"
$ws.Range("E5").Value = "2025-04-11 13:21:38"
$ws.Range("E6").Value = "2025-04-11 13:21:39"
$ws.Range("E7").Value = "2025-04-11 13:21:43"
$ws.Range("E8").Value = "2025-04-11 13:21:51"
$ws.Range("E9").Value = "2025-04-11 13:21:49"
$ws.Range("E10").Value = "2025-04-11 13:21:46"
$ws.Range("E11").Value = "2025-04-11 13:21:46"
$ws.Range("E12").Value = "2025-04-11 13:21:47"
$ws.Range("E13").Value = "2025-04-11 13:21:23"

# Correct the box id for row 13 (wrong scan corrected)
$ws.Range("B13").Value = 26002680643
